$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right
$ws.Columns.Item(1).Insert()

# Set the header for the newly inserted column A
$ws.Range("A1").Value = "ID"

# Update the selection to match the target state
$ws.Range("D6").Select()
